# Apply cell value updates across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the upstream diff (scheduled market-price data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4299.5264
$ws.Range("I51").Value = 2474.75
$ws.Range("J51").Value = 4786.1333
$ws.Range("K51").Value = 2474.75
$ws.Range("L51").Value = 4786.1333
$ws.Range("M51").Value = -1990.75
$ws.Range("N51").Value = -5754.1333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2482511
$ws.Range("I62").Value = 3482288.8
$ws.Range("J62").Value = 21519.54
$ws.Range("K62").Value = 3482288.8
$ws.Range("L62").Value = 21519.54
$ws.Range("M62").Value = -3481664.8
$ws.Range("N62").Value = -22767.54

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2482511
$ws.Range("I65").Value = 3482288.8
$ws.Range("J65").Value = 21519.54
$ws.Range("K65").Value = 17411444
$ws.Range("L65").Value = 107597.7
$ws.Range("M65").Value = -17408324
$ws.Range("N65").Value = -113837.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1490.0555
$ws.Range("I70").Value = 990.7273
$ws.Range("J70").Value = 2274.7144
$ws.Range("K70").Value = 2972.1819
$ws.Range("L70").Value = 6824.1432
$ws.Range("M70").Value = -2702.1819
$ws.Range("N70").Value = -7364.1432

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1490.0555
$ws.Range("I73").Value = 990.7273
$ws.Range("J73").Value = 2274.7144
$ws.Range("K73").Value = 2972.1819
$ws.Range("L73").Value = 6824.1432
$ws.Range("M73").Value = -2036.1819
$ws.Range("N73").Value = -8696.143199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3011.1667
$ws.Range("I116").Value = 3246.2307
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 3246.2307
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = 195.7692999999999
$ws.Range("N116").Value = -9284

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 22623.291
$ws.Range("I132").Value = 24425.408
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 73276.224
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -70746.224
$ws.Range("N132").Value = -13460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 37000
$ws.Range("J7").Value = 37000
$ws.Range("L7").Value = 37000
$ws.Range("N7").Value = -37228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 29599.5
$ws.Range("J114").Value = 29599.5
$ws.Range("L114").Value = 29599.5
$ws.Range("N114").Value = -38277.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 45671.5
$ws.Range("J135").Value = 45671.5
$ws.Range("L135").Value = 45671.5
$ws.Range("N135").Value = -55811.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 37293.332
$ws.Range("J81").Value = 37293.332
$ws.Range("L81").Value = 37293.332
$ws.Range("N81").Value = -39415.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 37293.332
$ws.Range("J84").Value = 37293.332
$ws.Range("L84").Value = 111879.996
$ws.Range("N84").Value = -122487.996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 53260
$ws.Range("J135").Value = 53260
$ws.Range("L135").Value = 53260
$ws.Range("N135").Value = -63400

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1952.5
$ws.Range("I35").Value = 1952.5
$ws.Range("K35").Value = 1952.5
$ws.Range("M35").Value = -1658.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2387.7083
$ws.Range("I99").Value = 1636.7142
$ws.Range("J99").Value = 3439.1
$ws.Range("K99").Value = 1636.7142
$ws.Range("L99").Value = 3439.1
$ws.Range("M99").Value = -138.7141999999999
$ws.Range("N99").Value = -6435.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 10371.75
$ws.Range("I103").Value = 4710.5713
$ws.Range("J103").Value = 50000
$ws.Range("K103").Value = 4710.5713
$ws.Range("L103").Value = 50000
$ws.Range("M103").Value = -3538.5713
$ws.Range("N103").Value = -52344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2387.7083
$ws.Range("I126").Value = 1636.7142
$ws.Range("J126").Value = 3439.1
$ws.Range("K126").Value = 4910.142599999999
$ws.Range("L126").Value = 10317.3
$ws.Range("M126").Value = -2440.142599999999
$ws.Range("N126").Value = -15257.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 882.9048
$ws.Range("I46").Value = 305
$ws.Range("J46").Value = 1063.5
$ws.Range("K46").Value = 915
$ws.Range("L46").Value = 3190.5
$ws.Range("M46").Value = -824
$ws.Range("N46").Value = -3372.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 35714930
$ws.Range("I113").Value = 637.5
$ws.Range("J113").Value = 50000644
$ws.Range("K113").Value = 1912.5
$ws.Range("L113").Value = 150001932
$ws.Range("M113").Value = 257.5
$ws.Range("N113").Value = -150006272

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 38669.332
$ws.Range("I22").Value = 38669.332
$ws.Range("K22").Value = 38669.332
$ws.Range("M22").Value = -38140.332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5296.324
$ws.Range("I70").Value = 5339.3
$ws.Range("K70").Value = 5339.3
$ws.Range("M70").Value = -5069.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5296.324
$ws.Range("I73").Value = 5339.3
$ws.Range("K73").Value = 5339.3
$ws.Range("M73").Value = -4403.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2747.72
$ws.Range("I102").Value = 2626
$ws.Range("K102").Value = 2626
$ws.Range("M102").Value = -1004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2550.9333
$ws.Range("I126").Value = 2194.3157
$ws.Range("J126").Value = 2811.5386
$ws.Range("K126").Value = 6582.9471
$ws.Range("L126").Value = 8434.6158
$ws.Range("M126").Value = -4112.9471
$ws.Range("N126").Value = -13374.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3205.1667
$ws.Range("I61").Value = 3205.1667
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3205.1667
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3003.1667
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 27400
$ws.Range("J92").Value = 27400
$ws.Range("L92").Value = 27400
$ws.Range("N92").Value = -32392

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3205.1667
$ws.Range("I113").Value = 3205.1667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3205.1667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1035.1667
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 10815.454
$ws.Range("I74").Value = 10523
$ws.Range("J74").Value = 10925.125
$ws.Range("K74").Value = 10523
$ws.Range("L74").Value = 10925.125
$ws.Range("M74").Value = -9587
$ws.Range("N74").Value = -12797.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 10815.454
$ws.Range("I77").Value = 10523
$ws.Range("J77").Value = 10925.125
$ws.Range("K77").Value = 31569
$ws.Range("L77").Value = 32775.375
$ws.Range("M77").Value = -26889
$ws.Range("N77").Value = -42135.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 669.44446
$ws.Range("I113").Value = 691.17645
$ws.Range("J113").Value = 300
$ws.Range("K113").Value = 2073.52935
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = 96.47064999999975
$ws.Range("N113").Value = -5240
